# Replace the 25 division-problem cells in the practice-sheet table.
# Table.Cell(row, col).Range.Text assignment is used (rather than Find.Execute)
# because this host's Find engine matches anywhere in the document, not just
# within the invoking Range -- direct Range.Text assignment stays cell-scoped.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text -like "91÷6=*") {
    $cell.Range.Text = "18÷5="
} else {
    Write-Host "Unexpected content in cell (1,1):" $cell.Range.Text
}
$cell = $t.Cell(1, 2)
if ($cell.Range.Text -like "82÷5=*") {
    $cell.Range.Text = "59÷2="
} else {
    Write-Host "Unexpected content in cell (1,2):" $cell.Range.Text
}
$cell = $t.Cell(1, 3)
if ($cell.Range.Text -like "55÷3=*") {
    $cell.Range.Text = "37÷9="
} else {
    Write-Host "Unexpected content in cell (1,3):" $cell.Range.Text
}
$cell = $t.Cell(1, 4)
if ($cell.Range.Text -like "44÷6=*") {
    $cell.Range.Text = "82÷8="
} else {
    Write-Host "Unexpected content in cell (1,4):" $cell.Range.Text
}
$cell = $t.Cell(1, 5)
if ($cell.Range.Text -like "80÷2=*") {
    $cell.Range.Text = "41÷6="
} else {
    Write-Host "Unexpected content in cell (1,5):" $cell.Range.Text
}
$cell = $t.Cell(5, 1)
if ($cell.Range.Text -like "16÷5=*") {
    $cell.Range.Text = "11÷4="
} else {
    Write-Host "Unexpected content in cell (5,1):" $cell.Range.Text
}
$cell = $t.Cell(5, 2)
if ($cell.Range.Text -like "25÷3=*") {
    $cell.Range.Text = "99÷3="
} else {
    Write-Host "Unexpected content in cell (5,2):" $cell.Range.Text
}
$cell = $t.Cell(5, 3)
if ($cell.Range.Text -like "41÷6=*") {
    $cell.Range.Text = "85÷2="
} else {
    Write-Host "Unexpected content in cell (5,3):" $cell.Range.Text
}
$cell = $t.Cell(5, 4)
if ($cell.Range.Text -like "20÷8=*") {
    $cell.Range.Text = "96÷7="
} else {
    Write-Host "Unexpected content in cell (5,4):" $cell.Range.Text
}
$cell = $t.Cell(5, 5)
if ($cell.Range.Text -like "39÷3=*") {
    $cell.Range.Text = "71÷6="
} else {
    Write-Host "Unexpected content in cell (5,5):" $cell.Range.Text
}
$cell = $t.Cell(9, 1)
if ($cell.Range.Text -like "54÷2=*") {
    $cell.Range.Text = "27÷9="
} else {
    Write-Host "Unexpected content in cell (9,1):" $cell.Range.Text
}
$cell = $t.Cell(9, 2)
if ($cell.Range.Text -like "34÷7=*") {
    $cell.Range.Text = "90÷9="
} else {
    Write-Host "Unexpected content in cell (9,2):" $cell.Range.Text
}
$cell = $t.Cell(9, 3)
if ($cell.Range.Text -like "65÷6=*") {
    $cell.Range.Text = "64÷8="
} else {
    Write-Host "Unexpected content in cell (9,3):" $cell.Range.Text
}
$cell = $t.Cell(9, 4)
if ($cell.Range.Text -like "68÷4=*") {
    $cell.Range.Text = "70÷3="
} else {
    Write-Host "Unexpected content in cell (9,4):" $cell.Range.Text
}
$cell = $t.Cell(9, 5)
if ($cell.Range.Text -like "31÷4=*") {
    $cell.Range.Text = "84÷2="
} else {
    Write-Host "Unexpected content in cell (9,5):" $cell.Range.Text
}
$cell = $t.Cell(13, 1)
if ($cell.Range.Text -like "79÷5=*") {
    $cell.Range.Text = "50÷8="
} else {
    Write-Host "Unexpected content in cell (13,1):" $cell.Range.Text
}
$cell = $t.Cell(13, 2)
if ($cell.Range.Text -like "20÷9=*") {
    $cell.Range.Text = "33÷2="
} else {
    Write-Host "Unexpected content in cell (13,2):" $cell.Range.Text
}
$cell = $t.Cell(13, 3)
if ($cell.Range.Text -like "51÷3=*") {
    $cell.Range.Text = "59÷2="
} else {
    Write-Host "Unexpected content in cell (13,3):" $cell.Range.Text
}
$cell = $t.Cell(13, 4)
if ($cell.Range.Text -like "91÷7=*") {
    $cell.Range.Text = "71÷6="
} else {
    Write-Host "Unexpected content in cell (13,4):" $cell.Range.Text
}
$cell = $t.Cell(13, 5)
if ($cell.Range.Text -like "51÷7=*") {
    $cell.Range.Text = "62÷3="
} else {
    Write-Host "Unexpected content in cell (13,5):" $cell.Range.Text
}
$cell = $t.Cell(17, 1)
if ($cell.Range.Text -like "26÷7=*") {
    $cell.Range.Text = "80÷9="
} else {
    Write-Host "Unexpected content in cell (17,1):" $cell.Range.Text
}
$cell = $t.Cell(17, 2)
if ($cell.Range.Text -like "50÷7=*") {
    $cell.Range.Text = "27÷8="
} else {
    Write-Host "Unexpected content in cell (17,2):" $cell.Range.Text
}
$cell = $t.Cell(17, 3)
if ($cell.Range.Text -like "90÷3=*") {
    $cell.Range.Text = "19÷9="
} else {
    Write-Host "Unexpected content in cell (17,3):" $cell.Range.Text
}
$cell = $t.Cell(17, 4)
if ($cell.Range.Text -like "44÷4=*") {
    $cell.Range.Text = "35÷7="
} else {
    Write-Host "Unexpected content in cell (17,4):" $cell.Range.Text
}
$cell = $t.Cell(17, 5)
if ($cell.Range.Text -like "36÷4=*") {
    $cell.Range.Text = "10÷5="
} else {
    Write-Host "Unexpected content in cell (17,5):" $cell.Range.Text
}

Write-Host "Replacements complete"
